# Weekly fruit/veg data update: a new price observation is inserted in the
# middle of the "Espinaca" time series (row 543), pushing the existing
# rows 543:571 down to 544:572 and growing the used range to A1:R572.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 543 (shifts 543:571 -> 544:572,
# inherits formatting -- including the date NumberFormat on column D --
# from the surrounding rows, just like Excel's UI "Insert Row" does).
$ws.Rows.Item(543).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A543").Value = 6
$ws.Range("B543").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C543").Value = "Metropolitana"
$ws.Range("D543").Value = 44753
$ws.Range("E543").Value = 13
$ws.Range("F543").Value = 100112012
$ws.Range("G543").Value = "Espinaca"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 340
$ws.Range("K543").Value = 9500
$ws.Range("L543").Value = 10000
$ws.Range("M543").Value = 9721
$ws.Range("N543").Value = "`$/cuna 10 kilos"
$ws.Range("O543").Value = "Región Metropolitana"
$ws.Range("P543").Value = 972
$ws.Range("Q543").Value = 10
$ws.Range("R543").Value = "Hortaliza"
